$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update the salary formula in column H (rows 3-32): multiply by 1200000*1.1 instead of 1000000
for ($r = 3; $r -le 32; $r++) {
    $ws.Range("H$r").Formula = "=(ROUND((G$r/1.5)*1200000*1.1,-5))/1000000"
}

# Update the active selection cell as recorded in the sheet view
$ws.Range("J18").Select()
